# feat: Add player and class filters to character stats
#
# Updates the character roster on Sheet1:
#   - D6 (吴工 / 邀月) class corrected from 骑士 to 圣骑士
#   - Rows 8:18 are re-sorted/grouped by player and several new characters
#     are inserted, growing the table down to row 27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Small, in-place fix: 吴工's "邀月" is a Paladin, not a generic Knight.
$ws.Range("D6").Value = "圣骑士"

# The rest of the roster (originally A8:D18) gets replaced wholesale with
# the regrouped-by-player list (now A8:D27), picking up several new
# characters along the way.
$rows = @(
    @('屯狗','屯屯宝宝','斯坦索姆','猎人'),
    @('巨奶','傻瓜观测','影之哀伤','牧师'),
    @('巨奶','天赐祝福','米奈希尔','圣骑士'),
    @('捷教授','四个自信','回音山','法师'),
    @('昶狂','冲锋先看路','霜之哀伤','战士'),
    @('昶狂','战复慢点起','霜之哀伤','死亡骑士'),
    @('段总','生锈的斩牛刀','伊森利恩','盗贼'),
    @('段总','飞翔的潼瑜','伊森利恩','死亡骑士'),
    @('统皇','焦糖扁可颂','斯坦索姆','圣骑士'),
    @('统皇','本间芽衣芓','斯坦索姆','战士'),
    @('统皇','生命众筹','斯坦索姆','死亡骑士'),
    @('统皇','亻沈默','图拉扬','法师'),
    @('舒总','Fountine','图拉扬','法师'),
    @('舒总','天灵浴血','诺兹多姆','死亡骑士'),
    @('舒总','霜满天丶','图拉扬','恶魔猎手'),
    @('蔡圣','莱恩弗尔特','神圣之歌','猎人'),
    @('蔡圣','亚妮艾丝','神圣之歌','牧师'),
    @('蔡圣','亚里欧斯','神圣之歌','恶魔猎手'),
    @('蔡圣','伊格瑞特','神圣之歌','德鲁伊'),
    @('蔡圣','萨里西翁','神圣之歌','萨满')
)

$startRow = 8
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Range("A$r").Value = $data[0]
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("C$r").Value = $data[2]
    $ws.Range("D$r").Value = $data[3]
}
